# Applies the "Add IQR var for SCE" edit described by the target diff.
#
# Summary of changes performed here (everything reachable through the
# Excel COM object model exposed by this runtime):
#   1. SPFPop (sheet2): update several regression-result text cells,
#      widen two bestFit-ish columns, move the active selection.
#   2. SPFInd (sheet1): move the active selection.
#   3. SCEInd (sheet4): move the active selection / scroll position.
#   4. Combined (sheet5): add a tab color, populate a small summary
#      table (this creates the new shared strings), and make it the
#      final active sheet (which also updates the workbook-level
#      activeTab and clears tabSelected on the previously active sheet).
#   5. Workbook: rename the FEEfficiencySPFQ defined name to
#      FEEfficiencySPFQ_1.
#
# NOTE: xl/connections.xml and xl/queryTables/*.xml are not reachable
# through this COM surface (Workbook.Connections / Worksheet.QueryTables
# are always empty in this runtime), so those parts of the original
# diff cannot be reproduced from here.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Defined name: FEEfficiencySPFQ -> FEEfficiencySPFQ_1
# ---------------------------------------------------------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "SPFPop!FEEfficiencySPFQ") {
        $n.Name = "FEEfficiencySPFQ_1"
    }
}

# ---------------------------------------------------------------------
# 2. SPFPop sheet: updated numbers + column widths
# ---------------------------------------------------------------------
$spfPop = $wb.Worksheets.Item("SPFPop")

# Row 3 - variable names
$spfPop.Range("B3").Formula = '="SPFCPI_FE_const"'
$spfPop.Range("F3").Formula = '="SPFPCE_FE_const"'

# Row 5/6 - L4.InfExp_Mean coefficient + se
$spfPop.Range("C5").Formula = '="0.243**"'
$spfPop.Range("G5").Formula = '="1.450"'
$spfPop.Range("C6").Formula = '="(0.089)"'
$spfPop.Range("G6").Formula = '="(0.865)"'

# Row 8/9 - L4.InfExp_FE coefficient + se
$spfPop.Range("D8").Formula = '="-0.00469"'
$spfPop.Range("H8").Formula = '="-0.290"'
$spfPop.Range("D9").Formula = '="(0.100)"'
$spfPop.Range("H9").Formula = '="(0.158)"'

# Row 11/12 - L.InfExp_FE coefficient + se
$spfPop.Range("E11").Formula = '="0.945***"'
$spfPop.Range("I11").Formula = '="1.026***"'
$spfPop.Range("E12").Formula = '="(0.117)"'
$spfPop.Range("I12").Formula = '="(0.209)"'

# Row 14/15 - L2.InfExp_FE coefficient + se
$spfPop.Range("E14").Formula = '="-0.258*"'
$spfPop.Range("I14").Formula = '="-0.491*"'
$spfPop.Range("E15").Formula = '="(0.119)"'
$spfPop.Range("I15").Formula = '="(0.186)"'

# Row 17/18 - L3.InfExp_FE coefficient + se
$spfPop.Range("E17").Formula = '="-0.0256"'
$spfPop.Range("I17").Formula = '="-0.0494"'
$spfPop.Range("E18").Formula = '="(0.097)"'
$spfPop.Range("I18").Formula = '="(0.179)"'

# Row 20/21 - _cons coefficient + se
$spfPop.Range("B20").Formula = '="0.308**"'
$spfPop.Range("C20").Formula = '="-0.500"'
$spfPop.Range("D20").Formula = '="0.244*"'
$spfPop.Range("E20").Formula = '="0.0845"'
$spfPop.Range("F20").Formula = '="0.279"'
$spfPop.Range("G20").Formula = '="-2.240"'
$spfPop.Range("H20").Formula = '="0.480*"'
$spfPop.Range("I20").Formula = '="0.216"'

$spfPop.Range("B21").Formula = '="(0.097)"'
$spfPop.Range("C21").Formula = '="(0.282)"'
$spfPop.Range("D21").Formula = '="(0.106)"'
$spfPop.Range("E21").Formula = '="(0.067)"'
$spfPop.Range("F21").Formula = '="(0.157)"'
$spfPop.Range("G21").Formula = '="(1.499)"'
$spfPop.Range("H21").Formula = '="(0.179)"'
$spfPop.Range("I21").Formula = '="(0.117)"'

# Row 23 - N
$spfPop.Range("B23").Formula = '="147"'
$spfPop.Range("C23").Formula = '="143"'
$spfPop.Range("D23").Formula = '="143"'
$spfPop.Range("E23").Formula = '="144"'
$spfPop.Range("F23").Formula = '="45"'
$spfPop.Range("G23").Formula = '="41"'
$spfPop.Range("H23").Formula = '="41"'
$spfPop.Range("I23").Formula = '="42"'

# Row 24 - r2
$spfPop.Range("C24").Formula = '="0.0683"'
$spfPop.Range("D24").Formula = '="0.0000247"'
$spfPop.Range("E24").Formula = '="0.594"'
$spfPop.Range("G24").Formula = '="0.122"'
$spfPop.Range("H24").Formula = '="0.103"'
$spfPop.Range("I24").Formula = '="0.654"'

# Column widths grew because the new labels ("...FE_const") render wider.
$spfPop.Columns.Item(2).ColumnWidth = 14
$spfPop.Columns.Item(6).ColumnWidth = 14.5

# ---------------------------------------------------------------------
# 3. Selections on sheets that are not the final active sheet.
#    These must be applied before we activate "Combined" at the end,
#    otherwise they would incorrectly become the active tab.
# ---------------------------------------------------------------------
$spfInd = $wb.Worksheets.Item("SPFInd")
$spfInd.Range("B20").Select()

$spfPop.Range("B4").Select()

$sceInd = $wb.Worksheets.Item("SCEInd")
$sceInd.Activate()
$excel.ActiveWindow.ScrollRow = 49
$sceInd.Range("C67").Select()

# ---------------------------------------------------------------------
# 4. Combined sheet: tab color + new summary table content.
#    Order of the Value assignments matters because it determines the
#    order in which new shared strings get interned.
# ---------------------------------------------------------------------
$combined = $wb.Worksheets.Item("Combined")

$combined.Range("A2").Value = "Test 1: Bias"
$combined.Range("B1").Value = "SPF CPI"
$combined.Range("C1").Value = "SPF PCE"
$combined.Range("D1").Value = "SCE"
$combined.Range("A7").Value = "Test2: FE Depends on past information"
$combined.Range("A12").Value = "Test3: FE of non-overllaping forecast horizons are serially correlated "
$combined.Range("A16").Value = "Test4: Overlapping FE are serially correlated "

$combined.Tab.Color = 65535
$combined.EnableFormatConditionsCalculation = $false

# Make "Combined" the final active sheet/selection so it ends up with
# tabSelected="1" and the workbook's activeTab points at it.
$combined.Activate()
$combined.Range("F23").Select()
